$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.152499999999996
$ws.Range("A12").Value = -21.6025
$ws.Range("D14").Value = -7.907599999999998
$ws.Range("D26").Value = -8.481600000000007
$ws.Range("D31").Value = -8.452499999999999
$ws.Range("A32").Value = -21.2258
$ws.Range("D35").Value = -8.357999999999997
$ws.Range("A36").Value = -20.1245
$ws.Range("D37").Value = -7.873599999999997
$ws.Range("A38").Value = -19.6076
$ws.Range("D45").Value = -7.723299999999998
$ws.Range("A46").Value = -21.8061
$ws.Range("A54").Value = -21.87439999999999
$ws.Range("A55").Value = -22.39810000000001
$ws.Range("D57").Value = -8.4626
$ws.Range("A67").Value = -21.41619999999997
$ws.Range("A69").Value = -21.56239999999997
$ws.Range("A72").Value = -21.93040000000001
$ws.Range("A91").Value = -21.4647
$ws.Range("A99").Value = -20.34839999999999
$ws.Range("D100").Value = -8.079099999999999
$ws.Range("D102").Value = -7.932799999999999
